$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.942.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.30%  "
$ws.Range("D3").Value = "'2.367.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'321.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.91%  "
$ws.Range("D6").Value = "'107.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.84%  "
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").Value = "'41.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.46%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "
$ws.Range("D12").Value = "'8.56"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.96%  "
$ws.Range("D13").Value = "'1.01"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.35%  "
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("D15").Value = "'16.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.57%  "
$ws.Range("D16").Value = "'2.726.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("D17").Value = "'2.426.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "'42.934.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.37%  "
$ws.Range("D19").Value = "'7.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'76.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").Value = "'3.73"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.80%  "
$ws.Range("D23").Value = "'261.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.71%  "
$ws.Range("D24").Value = "'2.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "
$ws.Range("D25").Value = "'9.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.75%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'11.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").Value = "'23.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.64%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.62%  "
$ws.Range("D30").Value = "'37.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'171.88"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.74%  "
$ws.Range("D32").Value = "'0.0903"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("E33").Value = "  -6.69%  "
$ws.Range("D34").Value = "'6.04"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").Value = "'0.122"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.79%  "
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").Value = "'4.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("E38").Value = "  +0.14%  "
$ws.Range("D39").Value = "'3.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.68%  "
$ws.Range("D40").Value = "'2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("D41").Value = "'1.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "'0.241"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.12%  "
$ws.Range("D43").Value = "'71.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("B44").Value = "BitcoinSV"
$ws.Range("C44").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D44").Value = "'96.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.64%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "'12.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.92%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'113.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.80%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D48").Value = "'5.57"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.41%  "
$ws.Range("D49").Value = "'9.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.33%  "
$ws.Range("D50").Value = "'77.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.77%  "
$ws.Range("D51").Value = "'1.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.11%  "
